$wb = $excel.ActiveWorkbook

# --- Update the time_taken values on the existing "data" sheet ---
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("F2").Value = "2021-10-05 14:35:35.554088"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:35.554096"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:35.554099"
$dataSheet.Range("F5").Value = "2021-10-05 14:35:35.554102"

# --- Add the new "metadata" sheet right after the "data" sheet ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# match page margins used by the "data" sheet
$metaSheet.PageSetup.LeftMargin = $dataSheet.PageSetup.LeftMargin
$metaSheet.PageSetup.RightMargin = $dataSheet.PageSetup.RightMargin
$metaSheet.PageSetup.TopMargin = $dataSheet.PageSetup.TopMargin
$metaSheet.PageSetup.BottomMargin = $dataSheet.PageSetup.BottomMargin
$metaSheet.PageSetup.HeaderMargin = $dataSheet.PageSetup.HeaderMargin
$metaSheet.PageSetup.FooterMargin = $dataSheet.PageSetup.FooterMargin

# Copy the header formatting (bold text, thin border, centered) from the
# data sheet's header row onto the new header row B1:G1 (and the A2 index
# cell), reusing the same cell style rather than creating a new one.
# -4122 = xlPasteFormats (paste formats only, leave cell contents alone).
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Renal Tubular Dysgenesis"
$metaSheet.Range("C2").Value = 166
$metaSheet.Range("D2").Value = "'0.4"
$metaSheet.Range("E2").Value = "2020-07-22T21:33:45.735788Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:35.550286"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/166/?format=json"

# D2 was entered with a leading quote to keep "0.4" as text rather than a
# number; re-apply the sheet's plain/unstyled format over it (content is
# untouched by a formats-only paste) so it doesn't carry a stray
# quote-prefix style
$dataSheet.Range("A1").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Ensure the "data" sheet stays active/selected
$dataSheet.Activate()
